$wb = $excel.ActiveWorkbook

# The existing "总计" sheet (index 6) already carries the right page setup /
# base styles, so clone it first -> the clone becomes the new "总计" roll-up
# sheet (placed last, picking up sheetId 7) while the original is repurposed
# below into the new "2022-Q1" per-fund holdings sheet (keeps sheetId 6).
$oldTotal = $wb.Worksheets.Item(6)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$oldTotal.Copy($null, $lastSheet)

$q1 = $wb.Worksheets.Item(6)
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Item(7)
$total.Name = "总计"

# ---------------------------------------------------------------------
# Sheet "2022-Q1": per-fund holding detail for the new quarter.
# ---------------------------------------------------------------------

# Headers (row 1). B:D already carry header style s=2 from the old sheet;
# clone that formatting onto the brand-new E:H header cells first.
$q1.Cells.Item(1,2).Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Column A (row index) keeps style s=2 for rows that already existed (2-6);
# row 7 is new, so clone the style from A2 onto it.
$q1.Cells.Item(2,1).Copy()
$q1.Cells.Item(7,1).PasteSpecial(-4122)

# Fund code (B) and the numeric-looking measures (D:G) must be stored as
# TEXT (fund codes carry significant leading zeros), matching the source.
$q1.Range("B2:B7").NumberFormat = "@"
$q1.Range("D2:G7").NumberFormat = "@"

$q1.Cells.Item(2,1).Value = 0
$q1.Cells.Item(2,2).Value = "169103"
$q1.Cells.Item(2,3).Value = "东方红睿轩三年定期开放灵活配置混合"
$q1.Cells.Item(2,4).Value = "40.03"
$q1.Cells.Item(2,5).Value = "57.34"
$q1.Cells.Item(2,6).Value = "2.69"
$q1.Cells.Item(2,7).Value = "1.0768"
$q1.Cells.Item(2,8).Value = 7

$q1.Cells.Item(3,1).Value = 1
$q1.Cells.Item(3,2).Value = "004278"
$q1.Cells.Item(3,3).Value = "东方红智逸沪港深定期开放混合"
$q1.Cells.Item(3,4).Value = "34.95"
$q1.Cells.Item(3,5).Value = "22.66"
$q1.Cells.Item(3,6).Value = "1.90"
$q1.Cells.Item(3,7).Value = "0.6640"
$q1.Cells.Item(3,8).Value = 5

$q1.Cells.Item(4,1).Value = 2
$q1.Cells.Item(4,2).Value = "012744"
$q1.Cells.Item(4,3).Value = "光大保德信品质生活混合型证券投资基金A"
$q1.Cells.Item(4,4).Value = "6.91"
$q1.Cells.Item(4,5).Value = "84.96"
$q1.Cells.Item(4,6).Value = "7.39"
$q1.Cells.Item(4,7).Value = "0.5106"
$q1.Cells.Item(4,8).Value = 2

$q1.Cells.Item(5,1).Value = 3
$q1.Cells.Item(5,2).Value = "003396"
$q1.Cells.Item(5,3).Value = "东方红优享红利沪港深灵活配置混合"
$q1.Cells.Item(5,4).Value = "15.87"
$q1.Cells.Item(5,5).Value = "89.27"
$q1.Cells.Item(5,6).Value = "3.01"
$q1.Cells.Item(5,7).Value = "0.4777"
$q1.Cells.Item(5,8).Value = 10

$q1.Cells.Item(6,1).Value = 4
$q1.Cells.Item(6,2).Value = "004099"
$q1.Cells.Item(6,3).Value = "前海开源沪港深景气行业精选灵活配置混合"
$q1.Cells.Item(6,4).Value = "0.41"
$q1.Cells.Item(6,5).Value = "93.07"
$q1.Cells.Item(6,6).Value = "8.87"
$q1.Cells.Item(6,7).Value = "0.0364"
$q1.Cells.Item(6,8).Value = 4

$q1.Cells.Item(7,1).Value = 5
$q1.Cells.Item(7,2).Value = "012758"
$q1.Cells.Item(7,3).Value = "光大保德信品质生活混合型证券投资基金C"
$q1.Cells.Item(7,4).Value = "0.31"
$q1.Cells.Item(7,5).Value = "84.96"
$q1.Cells.Item(7,6).Value = "7.39"
$q1.Cells.Item(7,7).Value = "0.0229"
$q1.Cells.Item(7,8).Value = 2

# ---------------------------------------------------------------------
# Sheet "总计": quarter-over-quarter roll-up, with 2022-Q1 prepended.
# ---------------------------------------------------------------------

# Row 7 is new (the table grew from 6 to 7 rows); clone column-A's index
# style from A2 so it matches the rest of the index column.
$total.Cells.Item(2,1).Copy()
$total.Cells.Item(7,1).PasteSpecial(-4122)

$total.Cells.Item(1,2).Value = "日期"
$total.Cells.Item(1,3).Value = "持有数量(只)"
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 6
$total.Cells.Item(2,4).Value = 2.79

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2021-Q4"
$total.Cells.Item(3,3).Value = 4
$total.Cells.Item(3,4).Value = 2.5

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2021-Q3"
$total.Cells.Item(4,3).Value = 6
$total.Cells.Item(4,4).Value = 5.59

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2021-Q2"
$total.Cells.Item(5,3).Value = 6
$total.Cells.Item(5,4).Value = 9.72

$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(6,2).Value = "2021-Q1"
$total.Cells.Item(6,3).Value = 9
$total.Cells.Item(6,4).Value = 11.1

$total.Cells.Item(7,1).Value = 5
$total.Cells.Item(7,2).Value = "2020-Q4"
$total.Cells.Item(7,3).Value = 9
$total.Cells.Item(7,4).Value = 9.19

# Restore the original active sheet/tab so the workbook view is undisturbed.
$wb.Worksheets.Item(1).Activate()
